# Regenerate merged AHB files
# - Rename the "_old" header columns (A:J) to "_FV2210"
# - Rename the "_new" header columns (L:U) to "_FV2304"
# - Freeze the header row (row 1)
# - Turn the A1:U88 range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) ---------------------------------------

$leftHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $leftHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}

# column K ("diff") stays the same

$rightHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $rightHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# --- 2. Freeze the header row ----------------------------------------------

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into a Table ------------------------------------

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U88"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
